$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell "pixel_size_mm" in I1, copying the style of the
# existing header cells (bold font) used for columns E:H (style index 2).
$ws.Cells.Item(1, 9).Value = "pixel_size_mm"
$ws.Cells.Item(1, 9).Font.Bold = $true

# Add the new data cell with the pixel size value.
$ws.Cells.Item(2, 9).Value = 1.8180000000000001

# Update the selection to match the post-edit state (I8).
$ws.Range("I8").Select()
